$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "distrOnly_rew"
$ws.Range("B1").Value = "distrOnly_potent"
$ws.Range("C1").Value = "distrOnly_percDead"

$ws.Range("A2").Value = "'0.6111111111111112"
$ws.Range("B2").Value = "'0.8777777777777778"
$ws.Range("C2").Value = 0.5

$ws.Range("A3").Value = "'0.6413043478260869"
$ws.Range("B3").Value = "'0.6630434782608695"
$ws.Range("C3").Value = "'0.16666666666666666"

$ws.Range("A4").Value = "'0.4298245614035088"
$ws.Range("B4").Value = "'0.5789473684210527"
$ws.Range("C4").Value = "'0.3333333333333333"

$ws.Range("A5").Value = "'0.5679012345679012"
$ws.Range("B5").Value = "'0.5679012345679012"
$ws.Range("C5").Value = 0.5

$ws.Range("A6").Value = "'0.7521367521367521"
$ws.Range("B6").Value = "'0.7521367521367521"
$ws.Range("C6").Value = "'0.16666666666666666"

$ws.Range("A7").Value = "'0.5161290322580645"
$ws.Range("B7").Value = "'0.6559139784946236"
$ws.Range("C7").Value = "'0.3333333333333333"

$ws.Range("A8").Value = "'0.8508771929824561"
$ws.Range("B8").Value = "'0.8947368421052632"
$ws.Range("C8").Value = "'0.16666666666666666"

$ws.Range("A9").Value = "'0.37777777777777777"
$ws.Range("B9").Value = "'0.6888888888888889"
$ws.Range("C9").Value = "'0.6666666666666666"

$ws.Range("A10").Value = 0.625
$ws.Range("B10").Value = "'0.8269230769230769"
$ws.Range("C10").Value = 0.5

$ws.Range("A11").Value = "'0.5913978494623656"
$ws.Range("B11").Value = "'0.6881720430107527"
$ws.Range("C11").Value = "'0.16666666666666666"

$ws.Range("A12").Value = 0.63157894736842102
$ws.Range("B12").Value = 0.63157894736842102
$ws.Range("C12").Value = "'0.16666666666666666"

$ws.Range("A13").Value = 0.69607843137254899
$ws.Range("B13").Value = "'0.7941176470588235"
$ws.Range("C13").Value = "'0.3333333333333333"

$ws.Range("A14").Value = "'0.8461538461538461"
$ws.Range("B14").Value = "'0.8461538461538461"
$ws.Range("C14").Value = 0

$ws.Range("A15").Value = 0.41121495327102803
$ws.Range("B15").Value = "'0.6074766355140186"
$ws.Range("C15").Value = "'0.6666666666666666"

[void]$ws.Range("G4").Select()
